$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Fbln1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.07168100000000001
$ws.Range("H2").Value = 0.215043
$ws.Range("I2").Value = 0.0008527427651799389
$ws.Range("J2").Value = 0.0008527427651799389
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 117.044563
$ws.Range("N2").Value = 351.133689
$ws.Range("O2").Value = 0.3245365645427815
$ws.Range("P2").Value = 0.3245365645427815
$ws.Range("Q2").Value = 8.389871320403001
$ws.Range("R2").Value = 75.508841883627
$ws.Range("S2").Value = 0.0002767462074502092
$ws.Range("T2").Value = 0.0002767462074502092
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Fbln1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.07168100000000001
$ws.Range("H3").Value = 0.215043
$ws.Range("I3").Value = 0.0008527427651799389
$ws.Range("J3").Value = 0.0008527427651799389
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 101.5800373333333
$ws.Range("N3").Value = 304.740112
$ws.Range("O3").Value = 0.281657135515876
$ws.Range("P3").Value = 0.281657135515876
$ws.Range("Q3").Value = 7.281358656090668
$ws.Range("R3").Value = 65.53222790481601
$ws.Range("S3").Value = 0.0002401810845724689
$ws.Range("T3").Value = 0.0002401810845724689
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Fbln1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.07168100000000001
$ws.Range("H4").Value = 0.215043
$ws.Range("I4").Value = 0.0008527427651799389
$ws.Range("J4").Value = 0.0008527427651799389
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 142.0267893333333
$ws.Range("N4").Value = 426.080368
$ws.Range("O4").Value = 0.3938062999413425
$ws.Range("P4").Value = 0.3938062999413425
$ws.Range("Q4").Value = 10.18062228620267
$ws.Range("R4").Value = 91.62560057582401
$ws.Range("S4").Value = 0.0003358154731572609
$ws.Range("T4").Value = 0.0003358154731572609
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Fbln1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 80.36585100000001
$ws.Range("H5").Value = 241.097553
$ws.Range("I5").Value = 0.9560608530542118
$ws.Range("J5").Value = 0.9560608530542118
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 117.044563
$ws.Range("N5").Value = 351.133689
$ws.Range("O5").Value = 0.3245365645427815
$ws.Range("P5").Value = 0.3245365645427815
$ws.Range("Q5").Value = 9406.385910418114
$ws.Range("R5").Value = 84657.47319376301
$ws.Range("S5").Value = 0.3102767047440549
$ws.Range("T5").Value = 0.3102767047440549
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Fbln1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 80.36585100000001
$ws.Range("H6").Value = 241.097553
$ws.Range("I6").Value = 0.9560608530542118
$ws.Range("J6").Value = 0.9560608530542118
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 101.5800373333333
$ws.Range("N6").Value = 304.740112
$ws.Range("O6").Value = 0.281657135515876
$ws.Range("P6").Value = 0.281657135515876
$ws.Range("Q6").Value = 8163.566144905105
$ws.Range("R6").Value = 73472.09530414594
$ws.Range("S6").Value = 0.2692813612501141
$ws.Range("T6").Value = 0.2692813612501141
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Fbln1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 80.36585100000001
$ws.Range("H7").Value = 241.097553
$ws.Range("I7").Value = 0.9560608530542118
$ws.Range("J7").Value = 0.9560608530542118
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 142.0267893333333
$ws.Range("N7").Value = 426.080368
$ws.Range("O7").Value = 0.3938062999413425
$ws.Range("P7").Value = 0.3938062999413425
$ws.Range("Q7").Value = 11414.10378957106
$ws.Range("R7").Value = 102726.9341061395
$ws.Range("S7").Value = 0.3765027870600428
$ws.Range("T7").Value = 0.3765027870600428
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Fbln1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 3.621815
$ws.Range("H8").Value = 10.865445
$ws.Range("I8").Value = 0.04308640418060826
$ws.Range("J8").Value = 0.04308640418060826
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 117.044563
$ws.Range("N8").Value = 351.133689
$ws.Range("O8").Value = 0.3245365645427815
$ws.Range("P8").Value = 0.3245365645427815
$ws.Range("Q8").Value = 423.913753941845
$ws.Range("R8").Value = 3815.223785476605
$ws.Range("S8").Value = 0.01398311359127634
$ws.Range("T8").Value = 0.01398311359127634
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Fbln1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 3.621815
$ws.Range("H9").Value = 10.865445
$ws.Range("I9").Value = 0.04308640418060826
$ws.Range("J9").Value = 0.04308640418060826
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 101.5800373333333
$ws.Range("N9").Value = 304.740112
$ws.Range("O9").Value = 0.281657135515876
$ws.Range("P9").Value = 0.281657135515876
$ws.Range("Q9").Value = 367.9041029144267
$ws.Range("R9").Value = 3311.13692622984
$ws.Range("S9").Value = 0.01213559318118939
$ws.Range("T9").Value = 0.01213559318118938
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Fbln1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 3.621815
$ws.Range("H10").Value = 10.865445
$ws.Range("I10").Value = 0.04308640418060826
$ws.Range("J10").Value = 0.04308640418060826
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 142.0267893333333
$ws.Range("N10").Value = 426.080368
$ws.Range("O10").Value = 0.3938062999413425
$ws.Range("P10").Value = 0.3938062999413425
$ws.Range("Q10").Value = 514.3947560093068
$ws.Range("R10").Value = 4629.552804083761
$ws.Range("S10").Value = 0.01696769740814253
$ws.Range("T10").Value = 0.01696769740814253
